$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3439.8
$ws.Range("I29").Value = 525.5
$ws.Range("J29").Value = 5382.6665
$ws.Range("K29").Value = 1576.5
$ws.Range("L29").Value = 16147.9995
$ws.Range("M29").Value = -1295.5
$ws.Range("N29").Value = -16709.9995
$ws.Range("H31").Value = 290.375
$ws.Range("I31").Value = 296.14285
$ws.Range("K31").Value = 888.4285500000001
$ws.Range("M31").Value = -658.4285500000001
$ws.Range("H34").Value = 2877.6428
$ws.Range("I34").Value = 2301.4
$ws.Range("K34").Value = 2301.4
$ws.Range("M34").Value = -2098.4
$ws.Range("H36").Value = 2877.6428
$ws.Range("I36").Value = 2301.4
$ws.Range("K36").Value = 2301.4
$ws.Range("M36").Value = -1586.4
$ws.Range("H58").Value = 173.125
$ws.Range("I58").Value = 173.125
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 519.375
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H76").Value = 5510.4287
$ws.Range("I76").Value = 5526.909
$ws.Range("J76").Value = 5450
$ws.Range("K76").Value = 5526.909
$ws.Range("L76").Value = 5450
$ws.Range("M76").Value = -5211.909
$ws.Range("N76").Value = -6080
$ws.Range("H79").Value = 5510.4287
$ws.Range("I79").Value = 5526.909
$ws.Range("J79").Value = 5450
$ws.Range("K79").Value = 5526.909
$ws.Range("L79").Value = 5450
$ws.Range("M79").Value = -4434.909
$ws.Range("N79").Value = -7634
$ws.Range("H92").Value = 392.35294
$ws.Range("I92").Value = 472.5
$ws.Range("J92").Value = 18.333334
$ws.Range("K92").Value = 472.5
$ws.Range("L92").Value = 18.333334
$ws.Range("M92").Value = 775.5
$ws.Range("N92").Value = -2514.333334
$ws.Range("H112").Value = 1393.7587
$ws.Range("J112").Value = 1401.5186
$ws.Range("L112").Value = 4204.5558
$ws.Range("N112").Value = -6420.5558
$ws.Range("H132").Value = 2454.4167
$ws.Range("I132").Value = 2148
$ws.Range("K132").Value = 6444
$ws.Range("M132").Value = -3914
$ws.Range("H138").Value = 2560.257
$ws.Range("J138").Value = 2131.6086
$ws.Range("L138").Value = 6394.825800000001
$ws.Range("N138").Value = -16674.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2848.2222
$ws.Range("I21").Value = 1324
$ws.Range("J21").Value = 3610.3333
$ws.Range("K21").Value = 1324
$ws.Range("L21").Value = 3610.3333
$ws.Range("M21").Value = -950
$ws.Range("N21").Value = -4358.3333
$ws.Range("H40").Value = 23968.857
$ws.Range("J40").Value = 27359.4
$ws.Range("L40").Value = 27359.4
$ws.Range("N40").Value = -27711.4
$ws.Range("H45").Value = 5490.375
$ws.Range("I45").Value = 6761.4
$ws.Range("J45").Value = 3372
$ws.Range("K45").Value = 6761.4
$ws.Range("L45").Value = 3372
$ws.Range("M45").Value = -6384.4
$ws.Range("N45").Value = -4126
$ws.Range("H61").Value = 2985.1428
$ws.Range("I61").Value = 2985.1428
$ws.Range("K61").Value = 2985.1428
$ws.Range("M61").Value = -2773.1428
$ws.Range("H132").Value = 2905.9412
$ws.Range("I132").Value = 1200.0834
$ws.Range("K132").Value = 3600.2502
$ws.Range("M132").Value = -1070.2502
$ws.Range("H136").Value = 2985.1428
$ws.Range("I136").Value = 2985.1428
$ws.Range("K136").Value = 8955.428400000001
$ws.Range("M136").Value = -6405.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 29990
$ws.Range("J9").Value = 29990
$ws.Range("L9").Value = 29990
$ws.Range("N9").Value = -30326
$ws.Range("H32").Value = 15000
$ws.Range("J32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("N32").Value = -15768
$ws.Range("H33").Value = 1214
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 8791
$ws.Range("I36").Value = 4423.6665
$ws.Range("J36").Value = 34995
$ws.Range("K36").Value = 4423.6665
$ws.Range("L36").Value = 34995
$ws.Range("M36").Value = -3889.6665
$ws.Range("N36").Value = -36063
$ws.Range("H76").Value = 18542.334
$ws.Range("J76").Value = 20313.5
$ws.Range("L76").Value = 20313.5
$ws.Range("N76").Value = -20943.5
$ws.Range("H79").Value = 18542.334
$ws.Range("J79").Value = 20313.5
$ws.Range("L79").Value = 20313.5
$ws.Range("N79").Value = -22497.5
$ws.Range("H94").Value = 2617.2646
$ws.Range("I94").Value = 2427.16
$ws.Range("K94").Value = 2427.16
$ws.Range("M94").Value = -1976.16
$ws.Range("H105").Value = 17403
$ws.Range("I105").Value = 17403
$ws.Range("K105").Value = 17403
$ws.Range("M105").Value = -15656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32702.545
$ws.Range("I31").Value = 29191.459
$ws.Range("K31").Value = 29191.459
$ws.Range("M31").Value = -28896.459
$ws.Range("H34").Value = 32702.545
$ws.Range("I34").Value = 29191.459
$ws.Range("K34").Value = 29191.459
$ws.Range("M34").Value = -28989.459
$ws.Range("H99").Value = 27028.723
$ws.Range("I99").Value = 32097.25
$ws.Range("K99").Value = 32097.25
$ws.Range("M99").Value = -30599.25
$ws.Range("H126").Value = 27028.723
$ws.Range("I126").Value = 32097.25
$ws.Range("K126").Value = 96291.75
$ws.Range("M126").Value = -93821.75
$ws.Range("H132").Value = 2671.4814
$ws.Range("I132").Value = 2599.6191
$ws.Range("J132").Value = 2923
$ws.Range("K132").Value = 7798.8573
$ws.Range("L132").Value = 8769
$ws.Range("M132").Value = -5268.8573
$ws.Range("N132").Value = -13829
$ws.Range("H134").Value = 2101
$ws.Range("I134").Value = 1686.9286
$ws.Range("K134").Value = 5060.7858
$ws.Range("M134").Value = -2525.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 451.6111
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H135").Value = 451.6111
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H43").Value = 25558.344
$ws.Range("J43").Value = 34562.523
$ws.Range("L43").Value = 34562.523
$ws.Range("N43").Value = -34864.523
$ws.Range("H80").Value = 2476.111
$ws.Range("I80").Value = 2219.1177
$ws.Range("K80").Value = 2219.1177
$ws.Range("M80").Value = -1221.1177
$ws.Range("H83").Value = 2476.111
$ws.Range("I83").Value = 2219.1177
$ws.Range("K83").Value = 11095.5885
$ws.Range("M83").Value = -6103.588499999998
$ws.Range("H97").Value = 23053.434
$ws.Range("I97").Value = 26372.96
$ws.Range("K97").Value = 26372.96
$ws.Range("M97").Value = -25876.96
$ws.Range("H126").Value = 3565.5625
$ws.Range("I126").Value = 3354.0833
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 10062.2499
$ws.Range("L126").Value = 12600
$ws.Range("M126").Value = -7592.249899999999
$ws.Range("N126").Value = -17540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6766.067
$ws.Range("I122").Value = 7618.2383
$ws.Range("J122").Value = 4777.6665
$ws.Range("K122").Value = 22854.7149
$ws.Range("L122").Value = 14332.9995
$ws.Range("M122").Value = -20404.7149
$ws.Range("N122").Value = -19232.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 26690.334
$ws.Range("J88").Value = 14950
$ws.Range("L88").Value = 14950
$ws.Range("N88").Value = -15762
$ws.Range("H91").Value = 26690.334
$ws.Range("J91").Value = 14950
$ws.Range("L91").Value = 14950
$ws.Range("N91").Value = -17758
$ws.Range("H132").Value = 5127.593
$ws.Range("I132").Value = 5615.227
$ws.Range("K132").Value = 16845.681
$ws.Range("M132").Value = -14315.681
